$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to stay a text value even when it "looks" numeric
    # (Excel would otherwise auto-convert e.g. "219.26" into a float and
    # drop significant trailing zeros). The leading apostrophe forces
    # text entry; resetting the Style back to Normal afterwards clears
    # the quotePrefix formatting flag/style index that the apostrophe
    # entry leaves behind, so the cell keeps the workbook's default style.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "28.312.73"
$ws.Range("E2").Value = "  +3.88%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.729.62"
$ws.Range("E3").Value = "  +2.18%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.04%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "219.26"
$ws.Range("E5").Value = "  +1.18%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.523"
$ws.Range("E6").Value = "  -0.14%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - Solana
Set-TextValue $ws.Range("D8") "23.99"

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.267"
$ws.Range("E9").Value = "  +1.50%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.30%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.56%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D12") "1.976.98"
$ws.Range("E12").Value = "  +2.40%  "

# Row 13 - WrappedEther
Set-TextValue $ws.Range("D13") "1.743.73"
$ws.Range("E13").Value = "  +3.11%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.96%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.58%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "67.75"
$ws.Range("E16").Value = "  +0.29%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "28.326.94"
$ws.Range("E17").Value = "  +3.97%  "

# Row 18 - BitcoinCash
Set-TextValue $ws.Range("D18") "246.81"
$ws.Range("E18").Value = "  +3.67%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0753"
$ws.Range("E19").Value = "  +0.85%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "7.91"
$ws.Range("E20").Value = "  -3.06%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.10%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.46%  "

# Row 23 - Avalanche
Set-TextValue $ws.Range("D23") "9.69"
$ws.Range("E23").Value = "  +0.16%  "

# Row 24 - Toncoin
Set-TextValue $ws.Range("D24") "2.09"
$ws.Range("E24").Value = "  -1.07%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "149.65"
$ws.Range("E25").Value = "  +0.87%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +2.24%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "16.70"
$ws.Range("E27").Value = "  +0.96%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +0.25%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.12%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +2.67%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +2.20%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.29%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "3.27"
$ws.Range("E33").Value = "  +0.53%  "

# Row 34 - Maker
Set-TextValue $ws.Range("D34") "1.491.59"
$ws.Range("E34").Value = "  -5.68%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -2.20%  "

# Row 36 - ARBITRUM
Set-TextValue $ws.Range("D36") "0.981"
$ws.Range("E36").Value = "  +2.40%  "

# Row 37 - ImmutableX
Set-TextValue $ws.Range("D37") "0.604"
$ws.Range("E37").Value = "  -0.28%  "

# Row 38 - HuobiToken
$ws.Range("E38").Value = "  +0.64%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +1.34%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  +0.27%  "

# Row 41 - Aave
Set-TextValue $ws.Range("D41") "70.03"
$ws.Range("E41").Value = "  +0.52%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.08%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  -2.02%  "

# Row 44 - MXToken
$ws.Range("E44").Value = "  +1.49%  "

# Row 45 - RocketPoolETH
Set-TextValue $ws.Range("D45") "1.880.62"
$ws.Range("E45").Value = "  +1.81%  "

# Row 46 - TrustWalletToken
Set-TextValue $ws.Range("D46") "0.801"
$ws.Range("E46").Value = "  +1.52%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  +6.97%  "

# Row 48 - now Quant (was BabyDogeCoin)
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D48") "90.63"
$ws.Range("E48").Value = "  -0.82%  "

# Row 49 - now BabyDogeCoin (was Quant)
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0113"
$ws.Range("E49").Value = "  +3.65%  "

# Row 50 - EnergySwap
$ws.Range("E50").Value = "  -1.00%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  -1.17%  "
